# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# zh-cn and de-de sheets (rows 2 and 4 share identical text values).

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-13 04:19:55"
$wsZh.Range("E4").Value = "2016-03-13 04:19:55"
$wsZh.Range("H2").Value = "2016-03-13 04:20:28"
$wsZh.Range("H4").Value = "2016-03-13 04:20:28"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-13 04:19:59"
$wsDe.Range("E4").Value = "2016-03-13 04:19:59"
$wsDe.Range("H2").Value = "2016-03-13 04:20:34"
$wsDe.Range("H4").Value = "2016-03-13 04:20:34"
